# UGCSyllabus.xlsx update: replace Unit-8 (Family, Marriage & Kinship) content
# with Unit-9 (Science, Technology & Society) content, and extend the
# column-A formatting (that used to only cover a handful of "block border"
# rows) down across the whole A17:A123 grid - mirroring the existing
# left/center/wrap (plain + bold) and general/center formats already used
# elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the unit heading + 15 topic rows (A1:A16) with the new content.
# ---------------------------------------------------------------------------

$ws.Range("A1").Value2 = "Unit - 9 : Science, Technology & Society`nTechnology and Emerging Political Processes, State Policy, Digital Divide and Inclusion"
$ws.Range("A2").Value2 = "Networked social movements and connective action (Castells; Bennett and Segerberg)"
$ws.Range("A3").Value2 = "Hashtag activism versus slacktivism"
$ws.Range("A4").Value2 = "Computational propaganda, bots, and astroturfing (Philip N. Howard)"
$ws.Range("A5").Value2 = "Microtargeting and political advertising on platforms"
$ws.Range("A6").Value2 = "E-participation tools: e-petitions, MyGov, online RTI"
$ws.Range("A7").Value2 = "Civic technology and open data for accountability"
$ws.Range("A8").Value2 = "Platform governance and content moderation as political process"
$ws.Range("A9").Value2 = "Internet shutdowns and throttling as digital authoritarian practices"
$ws.Range("A10").Value2 = "Data localization and digital sovereignty"
$ws.Range("A11").Value2 = "Digital divide levels: access, skills, usage, outcome divides"
$ws.Range("A12").Value2 = "Intersectional divides: gender, rural–urban, disability, language"
$ws.Range("A13").Value2 = "Inclusion strategies: BharatNet, CSCs, PMGDISHA, WCAG accessibility"
$ws.Range("A14").Value2 = "Assistive technology and inclusive design: screen readers, captions, UPI 123PAY"
$ws.Range("A15").Value2 = "EVMs and VVPAT: trust, transparency, auditability debates"
$ws.Range("A16").Value2 = "Platform/gig workers’ collective action and algorithmic bargaining"

# Entering the multi-line title re-triggers row autofit in row 1; pin the
# row height back to the sheet's standard 13.95 so the rest of the grid
# stays visually consistent.
$ws.Rows(1).RowHeight = 13.95

# ---------------------------------------------------------------------------
# 2. Extend formatting down column A for the (previously mostly blank)
#    grid rows 17:123, reusing the formats already present elsewhere on
#    the sheet (left/center/wrap "plain", left/center/wrap "bold", and
#    general/center "no-wrap").
# ---------------------------------------------------------------------------

function Set-PlainWrapLeft($addr) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4131   # xlLeft
    $r.VerticalAlignment = -4108     # xlCenter
    $r.WrapText = $true
    $r.Font.Bold = $false
}

function Set-BoldWrapLeft($addr) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4131   # xlLeft
    $r.VerticalAlignment = -4108     # xlCenter
    $r.WrapText = $true
    $r.Font.Bold = $true
}

function Set-GeneralNoWrap($addr) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = 1       # xlGeneral
    $r.VerticalAlignment = -4108     # xlCenter
    $r.WrapText = $false
    $r.Font.Bold = $false
}

# left/center/wrap, regular weight
Set-PlainWrapLeft("A17:A24")
Set-PlainWrapLeft("A26:A32")
Set-PlainWrapLeft("A65:A75")
Set-PlainWrapLeft("A77:A91")
Set-PlainWrapLeft("A93:A107")
Set-PlainWrapLeft("A109:A123")

# left/center/wrap, bold (the "section boundary" rows)
Set-BoldWrapLeft("A25")
Set-BoldWrapLeft("A41")
Set-BoldWrapLeft("A60")
Set-BoldWrapLeft("A76")
Set-BoldWrapLeft("A92")
Set-BoldWrapLeft("A108")

# general/center, no wrap
Set-GeneralNoWrap("A33:A40")
Set-GeneralNoWrap("A42:A56")
Set-GeneralNoWrap("A57:A59")
Set-GeneralNoWrap("A61:A64")

# ---------------------------------------------------------------------------
# 3. Selection: mirror the saved workbook's cursor state (A1:A123 selected).
# ---------------------------------------------------------------------------

$ws.Range("A1:A123").Select()
